$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 131.16667
$ws.Range("I9").Value = 129.75
$ws.Range("J9").Value = 134
$ws.Range("K9").Value = 129.75
$ws.Range("L9").Value = 134
$ws.Range("M9").Value = 39.25
$ws.Range("N9").Value = -472
$ws.Range("H12").Value = 500
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -330
$ws.Range("H17").Value = 3987.25
$ws.Range("J17").Value = 3987.25
$ws.Range("L17").Value = 11961.75
$ws.Range("N17").Value = -12297.75
$ws.Range("H42").Value = 152.8
$ws.Range("I42").Value = 24.5
$ws.Range("K42").Value = 73.5
$ws.Range("M42").Value = 156.5
$ws.Range("H113").Value = 13921.277
$ws.Range("I113").Value = 2863.625
$ws.Range("J113").Value = 22767.4
$ws.Range("K113").Value = 2863.625
$ws.Range("L113").Value = 22767.4
$ws.Range("M113").Value = 390.375
$ws.Range("N113").Value = -29275.4
$ws.Range("H138").Value = 6164.5425
$ws.Range("I138").Value = 7317.6816
$ws.Range("J138").Value = 5478.892
$ws.Range("K138").Value = 21953.0448
$ws.Range("L138").Value = 16436.676
$ws.Range("M138").Value = -16813.0448
$ws.Range("N138").Value = -26716.676

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22321.623
$ws.Range("I32").Value = 22237.5
$ws.Range("K32").Value = 22237.5
$ws.Range("M32").Value = -21950.5
$ws.Range("H45").Value = 3806.842
$ws.Range("I45").Value = 2828.8
$ws.Range("K45").Value = 2828.8
$ws.Range("M45").Value = -2451.8
$ws.Range("H55").Value = 37000
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H61").Value = 6143.727
$ws.Range("I61").Value = 6063.1
$ws.Range("J61").Value = 6950
$ws.Range("K61").Value = 6063.1
$ws.Range("L61").Value = 6950
$ws.Range("M61").Value = -5851.1
$ws.Range("N61").Value = -7374
$ws.Range("H97").Value = 754.0345
$ws.Range("I97").Value = 854.16
$ws.Range("K97").Value = 854.16
$ws.Range("M97").Value = -358.16
$ws.Range("H110").Value = 3295.4707
$ws.Range("I110").Value = 772.61536
$ws.Range("K110").Value = 772.61536
$ws.Range("M110").Value = 1272.38464
$ws.Range("H132").Value = 2294.7585
$ws.Range("I132").Value = 2277.18
$ws.Range("K132").Value = 6831.539999999999
$ws.Range("M132").Value = -4301.539999999999
$ws.Range("H136").Value = 6143.727
$ws.Range("I136").Value = 6063.1
$ws.Range("J136").Value = 6950
$ws.Range("K136").Value = 18189.3
$ws.Range("L136").Value = 20850
$ws.Range("M136").Value = -15639.3
$ws.Range("N136").Value = -25950

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 932.6667
$ws.Range("I11").Value = 800
$ws.Range("J11").Value = 999
$ws.Range("K11").Value = 800
$ws.Range("L11").Value = 999
$ws.Range("M11").Value = -660
$ws.Range("N11").Value = -1279
$ws.Range("H20").Value = 6492.5454
$ws.Range("J20").Value = 8149.75
$ws.Range("L20").Value = 8149.75
$ws.Range("N20").Value = -8643.75
$ws.Range("H82").Value = 10934.637
$ws.Range("I82").Value = 3364.5557
$ws.Range("J82").Value = 45000
$ws.Range("K82").Value = 3364.5557
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -2981.5557
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 10934.637
$ws.Range("I85").Value = 3364.5557
$ws.Range("J85").Value = 45000
$ws.Range("K85").Value = 3364.5557
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -2038.5557
$ws.Range("N85").Value = -47652
$ws.Range("H105").Value = 1691.25
$ws.Range("I105").Value = 1539.8334
$ws.Range("J105").Value = 2599.75
$ws.Range("K105").Value = 1539.8334
$ws.Range("L105").Value = 2599.75
$ws.Range("M105").Value = 207.1666
$ws.Range("N105").Value = -6093.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7062.077
$ws.Range("I31").Value = 5551.5
$ws.Range("K31").Value = 5551.5
$ws.Range("M31").Value = -5256.5
$ws.Range("H34").Value = 7062.077
$ws.Range("I34").Value = 5551.5
$ws.Range("K34").Value = 5551.5
$ws.Range("M34").Value = -5349.5
$ws.Range("H41").Value = 27333.334
$ws.Range("J41").Value = 23500
$ws.Range("L41").Value = 23500
$ws.Range("N41").Value = -24356
$ws.Range("H58").Value = 5261.2856
$ws.Range("I58").Value = 3192
$ws.Range("K58").Value = 3192
$ws.Range("M58").Value = -2989
$ws.Range("H62").Value = 6158.5
$ws.Range("J62").Value = 8289.666999999999
$ws.Range("L62").Value = 8289.666999999999
$ws.Range("N62").Value = -9537.666999999999
$ws.Range("H65").Value = 6158.5
$ws.Range("J65").Value = 8289.666999999999
$ws.Range("L65").Value = 41448.335
$ws.Range("N65").Value = -47688.335
$ws.Range("H70").Value = 42333.332
$ws.Range("J70").Value = 42333.332
$ws.Range("L70").Value = 42333.332
$ws.Range("N70").Value = -42963.332
$ws.Range("H73").Value = 42333.332
$ws.Range("J73").Value = 42333.332
$ws.Range("L73").Value = 42333.332
$ws.Range("N73").Value = -44517.332
$ws.Range("H99").Value = 8098.154
$ws.Range("J99").Value = 9247.125
$ws.Range("L99").Value = 9247.125
$ws.Range("N99").Value = -12243.125
$ws.Range("H126").Value = 8098.154
$ws.Range("J126").Value = 9247.125
$ws.Range("L126").Value = 27741.375
$ws.Range("N126").Value = -32681.375
$ws.Range("H134").Value = 2566.2
$ws.Range("I134").Value = 1533.3529
$ws.Range("K134").Value = 4600.0587
$ws.Range("M134").Value = -2065.0587
$ws.Range("H136").Value = 5261.2856
$ws.Range("I136").Value = 3192
$ws.Range("K136").Value = 9576
$ws.Range("M136").Value = -7026
$ws.Range("H141").Value = 299597.06
$ws.Range("J141").Value = 325784.94
$ws.Range("L141").Value = 325784.94
$ws.Range("N141").Value = -336144.94

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 8612
$ws.Range("I141").Value = 8612
$ws.Range("K141").Value = 25836
$ws.Range("M141").Value = -20656

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10550
$ws.Range("I70").Value = 8500
$ws.Range("K70").Value = 8500
$ws.Range("M70").Value = -8230
$ws.Range("H73").Value = 10550
$ws.Range("I73").Value = 8500
$ws.Range("K73").Value = 8500
$ws.Range("M73").Value = -7564
$ws.Range("H132").Value = 2537.85
$ws.Range("I132").Value = 2179.0588
$ws.Range("K132").Value = 6537.176399999999
$ws.Range("M132").Value = -4007.176399999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6503.56
$ws.Range("I7").Value = 5783.9473
$ws.Range("K7").Value = 5783.9473
$ws.Range("M7").Value = -5671.9473
$ws.Range("H16").Value = 482.5
$ws.Range("I16").Value = 440.66666
$ws.Range("J16").Value = 566.1667
$ws.Range("K16").Value = 440.66666
$ws.Range("L16").Value = 566.1667
$ws.Range("M16").Value = -270.66666
$ws.Range("N16").Value = -906.1667
$ws.Range("H55").Value = 466.66666
$ws.Range("I55").Value = 400
$ws.Range("K55").Value = 400
$ws.Range("M55").Value = -227
$ws.Range("H100").Value = 6811.385
$ws.Range("I100").Value = 4125.4165
$ws.Range("K100").Value = 4125.4165
$ws.Range("M100").Value = -3584.4165
$ws.Range("H126").Value = 6503.56
$ws.Range("I126").Value = 5783.9473
$ws.Range("K126").Value = 17351.8419
$ws.Range("M126").Value = -14881.8419

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 7354.75
$ws.Range("J15").Value = 7318.5
$ws.Range("L15").Value = 7318.5
$ws.Range("N15").Value = -7894.5
$ws.Range("H104").Value = 26789.334
$ws.Range("J104").Value = 26789.334
$ws.Range("L104").Value = 26789.334
$ws.Range("N104").Value = -33777.334
$ws.Range("H107").Value = 1282.3939
$ws.Range("I107").Value = 961.9167
$ws.Range("J107").Value = 2137
$ws.Range("K107").Value = 2885.7501
$ws.Range("L107").Value = 6411
$ws.Range("M107").Value = -965.7501000000002
$ws.Range("N107").Value = -10251
$ws.Range("H113").Value = 2064.04
$ws.Range("I113").Value = 1180.25
$ws.Range("J113").Value = 5599.2
$ws.Range("K113").Value = 3540.75
$ws.Range("L113").Value = 16797.6
$ws.Range("M113").Value = -1370.75
$ws.Range("N113").Value = -21137.6
$ws.Range("H132").Value = 131392.42
$ws.Range("I132").Value = 164004.66
$ws.Range("K132").Value = 492013.98
$ws.Range("M132").Value = -489483.98
$ws.Range("H136").Value = 3036.9575
$ws.Range("I136").Value = 1661.1818
$ws.Range("J136").Value = 6279.857
$ws.Range("K136").Value = 4983.5454
$ws.Range("L136").Value = 18839.571
$ws.Range("M136").Value = -2433.5454
$ws.Range("N136").Value = -23939.571
